$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Case Worker Data")
$ws.Name = "Staff Data"
Write-Output $wb.Worksheets.Item(3).Name
